$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text (it stores price strings like "60.362.79")
# while we update values, then restore default styling afterward so the
# underlying cell style (s attribute) is left untouched.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.362.79"
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").Value = "2.321.90"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "545.34"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "131.04"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.580"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").Value = "2.319.54"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").Value = "5.49"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "23.67"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "60.298.08"
$ws.Range("E15").Value = "  +3.42%  "
$ws.Range("D16").Value = "2.732.55"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "2.317.89"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("D19").Value = "10.59"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "4.14"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "313.75"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "6.62"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "63.80"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "7.84"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").Value = "1.34"
$ws.Range("E28").Value = "  +3.76%  "
$ws.Range("D29").Value = "173.55"
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.74"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("B31").Value = "SuiNetwork"
$ws.Range("C31").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").Value = "  +7.52%  "
$ws.Range("D32").Value = "0.0₃0730"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "5.92"
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("E34").Value = "  +10.28%  "
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "17.84"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "4.05"
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("D40").Value = "323.73"
$ws.Range("E40").Value = "  +11.42%  "
$ws.Range("D41").Value = "37.94"
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").Value = "138.03"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("D44").Value = "3.49"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").Value = "0.0941"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("D46").Value = "19.23"
$ws.Range("E46").Value = "  +4.98%  "
$ws.Range("D47").Value = "0.0496"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "0.559"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "0.0213"
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "11.03"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0210"
$ws.Range("E51").Value = "  +16.46%  "

# Restore original (default) style on column D now that values are set as text.
$ws.Range("D2:D51").Style = "Normal"
